# Weekly price-sheet update: a new weekly reading (Fecha = 2022-12-28,
# serial 44923) is inserted above the existing history, pushing every
# existing record down by two rows. Excel stores the reading twice
# (rows 36 and 37 end up identical), exactly as in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (previously rows 36:65) down by two rows so we
# can drop the new readings in at rows 36:37. Inserting at row 36 twice
# (rather than a single 2-row insert) matches Excel's default behaviour
# of carrying the above row's number format (column D's date style) into
# each freshly inserted row.
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(36).Insert()

$newRowsData = @(
    @{ Row = 36; Fecha = 44923; Volumen = 480; PMin = 1500; PMax = 1500; PProm = 1500; PKg = 1500 },
    @{ Row = 37; Fecha = 44923; Volumen = 480; PMin = 1500; PMax = 1500; PProm = 1500; PKg = 1500 }
)

foreach ($entry in $newRowsData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value  = 3
    $ws.Cells.Item($r, 2).Value  = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value  = "Coquimbo"
    $ws.Cells.Item($r, 4).Value  = $entry.Fecha
    $ws.Cells.Item($r, 5).Value  = 5
    $ws.Cells.Item($r, 6).Value  = 300000000
    $ws.Cells.Item($r, 7).Value  = "Espárragos"
    $ws.Cells.Item($r, 8).Value  = "Verde"
    $ws.Cells.Item($r, 9).Value  = "Primera"
    $ws.Cells.Item($r, 10).Value = $entry.Volumen
    $ws.Cells.Item($r, 11).Value = $entry.PMin
    $ws.Cells.Item($r, 12).Value = $entry.PMax
    $ws.Cells.Item($r, 13).Value = $entry.PProm
    $ws.Cells.Item($r, 14).Value = "$/kilo"
    $ws.Cells.Item($r, 15).Value = "Provincia de Quillota"
    $ws.Cells.Item($r, 16).Value = $entry.PKg
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
